$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Retanqueo")

$ws.Range("T2").Value = '"03/03/2022"'
$ws.Range("C2").Value = '"99509"'
$ws.Range("E2").Value = '"1.71"'
$ws.Range("N2").Value = '"Marzo"'
$ws.Range("B2").Value = '"7255282"'
$ws.Range("H2").Value = '"2036000"'
$ws.Range("J2").Value = '"871228"'
$ws.Range("M2").Value = '"JAIVER"'
$ws.Range("F2").Value = '"120"'
$ws.Range("I2").Value = '"0"'
$ws.Range("Z2").Value = '"03/03/2022"'

$ws.Range("M4").Select()
